# Add a new "Ocean" terrain column (F) to the Locations table on Sheet2.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# New header + first two entries for the Ocean terrain category.
$ws.Range("F1").Value = "Ocean"
$ws.Range("F2").Value = "Shipwreck"
$ws.Range("F3").Value = "Sucken Temple"

# Leave the selection where the author left it after typing the new column.
$ws.Range("F4").Select()
